# Fix year/area values (shift up by 9 rows) and round D column to 2 decimals,
# then remove rows 24-36 (years 2021-2024 entries that had no D value) so the
# sheet shrinks back down to A1:D23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Year / Area km^2 pairs for rows 2..23 (shifted up by 9 from the original table)
$years = @(1999,2000,2001,2002,2003,2004,2005,2006,2007,2008,2009,2010,2011,2012,2013,2014,2015,2016,2017,2018,2019,2020)
$areas = @(5111,6671,5237,7510,7145,8870,5899,5659,5526,5607,4281,3770,3008,1741,2346,1887,2153,2992,2433,2744,4172,4899)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $years[$i]
    $ws.Cells.Item($row, 2).Value = $areas[$i]
}

# Round the existing D column values (rows 2..23) to 2 decimal places in place
for ($row = 2; $row -le 23; $row++) {
    $cell = $ws.Cells.Item($row, 4)
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value = [Math]::Round([double]$current, 2)
    }
}

# Delete rows 24 through 36 (no longer part of the data set)
$ws.Rows("24:36").Delete()
